$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 241.75
$ws.Range("I2").Value = 241.75
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 241.75
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -128.75
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 40
$ws.Range("I4").Value = 40
$ws.Range("K4").Value = 40
$ws.Range("M4").Value = 74

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 844.9286
$ws.Range("I18").Value = 844.9286
$ws.Range("K18").Value = 844.9286
$ws.Range("M18").Value = -560.9286

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 631.6667
$ws.Range("I39").Value = 52.5
$ws.Range("J39").Value = 797.1429000000001
$ws.Range("K39").Value = 157.5
$ws.Range("L39").Value = 2391.4287
$ws.Range("M39").Value = 138.5
$ws.Range("N39").Value = -2983.4287

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2676.25
$ws.Range("I43").Value = 1850.5
$ws.Range("K43").Value = 1850.5
$ws.Range("M43").Value = -1781.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1036.9697
$ws.Range("I112").Value = 415
$ws.Range("J112").Value = 1236
$ws.Range("K112").Value = 1245
$ws.Range("L112").Value = 3708
$ws.Range("M112").Value = -137
$ws.Range("N112").Value = -5924

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3018.2551
$ws.Range("I137").Value = 2620.451
$ws.Range("K137").Value = 7861.353
$ws.Range("M137").Value = -5311.353

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3954.075
$ws.Range("I141").Value = 1488.6666
$ws.Range("K141").Value = 4465.9998
$ws.Range("M141").Value = 714.0002000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2698.6094
$ws.Range("I32").Value = 2615.8728
$ws.Range("K32").Value = 2615.8728
$ws.Range("M32").Value = -2328.8728

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 10500
$ws.Range("J56").Value = 10500
$ws.Range("L56").Value = 10500
$ws.Range("N56").Value = -11984

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 34484652
$ws.Range("I61").Value = 50001716
$ws.Range("J61").Value = 2289.7778
$ws.Range("K61").Value = 50001716
$ws.Range("L61").Value = 2289.7778
$ws.Range("M61").Value = -50001504
$ws.Range("N61").Value = -2713.7778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5397.7207
$ws.Range("I74").Value = 6468.4116
$ws.Range("J74").Value = 4697.654
$ws.Range("K74").Value = 6468.4116
$ws.Range("L74").Value = 4697.654
$ws.Range("M74").Value = -5594.4116
$ws.Range("N74").Value = -6445.654

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5397.7207
$ws.Range("I77").Value = 6468.4116
$ws.Range("J77").Value = 4697.654
$ws.Range("K77").Value = 32342.058
$ws.Range("L77").Value = 23488.27
$ws.Range("M77").Value = -27974.058
$ws.Range("N77").Value = -32224.27

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6255.56
$ws.Range("I132").Value = 5279.1763
$ws.Range("J132").Value = 8330.375
$ws.Range("K132").Value = 15837.5289
$ws.Range("L132").Value = 24991.125
$ws.Range("M132").Value = -13307.5289
$ws.Range("N132").Value = -30051.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 34484652
$ws.Range("I136").Value = 50001716
$ws.Range("J136").Value = 2289.7778
$ws.Range("K136").Value = 150005148
$ws.Range("L136").Value = 6869.3334
$ws.Range("M136").Value = -150002598
$ws.Range("N136").Value = -11969.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2225.04
$ws.Range("I134").Value = 1922.1111
$ws.Range("J134").Value = 3004
$ws.Range("K134").Value = 5766.3333
$ws.Range("L134").Value = 9012
$ws.Range("M134").Value = -3231.3333
$ws.Range("N134").Value = -14082

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4684.5317
$ws.Range("I31").Value = 1062.5
$ws.Range("J31").Value = 6461.3774
$ws.Range("K31").Value = 1062.5
$ws.Range("L31").Value = 6461.3774
$ws.Range("M31").Value = -767.5
$ws.Range("N31").Value = -7051.3774

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4684.5317
$ws.Range("I34").Value = 1062.5
$ws.Range("J34").Value = 6461.3774
$ws.Range("K34").Value = 1062.5
$ws.Range("L34").Value = 6461.3774
$ws.Range("M34").Value = -860.5
$ws.Range("N34").Value = -6865.3774

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 47626256
$ws.Range("I132").Value = 111124800
$ws.Range("J132").Value = 2344.5
$ws.Range("K132").Value = 333374400
$ws.Range("L132").Value = 7033.5
$ws.Range("M132").Value = -333371870
$ws.Range("N132").Value = -12093.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3594.9546
$ws.Range("I134").Value = 4038.8462
$ws.Range("K134").Value = 12116.5386
$ws.Range("M134").Value = -9581.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7789.2383
$ws.Range("I132").Value = 13391.556
$ws.Range("J132").Value = 3587.5
$ws.Range("K132").Value = 40174.66800000001
$ws.Range("L132").Value = 10762.5
$ws.Range("M132").Value = -37644.66800000001
$ws.Range("N132").Value = -15822.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2323.7368
$ws.Range("I22").Value = 2864.1428
$ws.Range("J22").Value = 2008.5
$ws.Range("K22").Value = 2864.1428
$ws.Range("L22").Value = 2008.5
$ws.Range("M22").Value = -2569.1428
$ws.Range("N22").Value = -2598.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2323.7368
$ws.Range("I27").Value = 2864.1428
$ws.Range("J27").Value = 2008.5
$ws.Range("K27").Value = 2864.1428
$ws.Range("L27").Value = 2008.5
$ws.Range("M27").Value = -2757.1428
$ws.Range("N27").Value = -2222.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2760.147
$ws.Range("I40").Value = 2578.1428
$ws.Range("J40").Value = 3054.1538
$ws.Range("K40").Value = 2578.1428
$ws.Range("L40").Value = 3054.1538
$ws.Range("M40").Value = -2442.1428
$ws.Range("N40").Value = -3326.1538

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5358.162
$ws.Range("I132").Value = 5839.3447
$ws.Range("J132").Value = 3613.875
$ws.Range("K132").Value = 17518.0341
$ws.Range("L132").Value = 10841.625
$ws.Range("M132").Value = -14988.0341
$ws.Range("N132").Value = -15901.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4209.2856
$ws.Range("I136").Value = 1113.5217
$ws.Range("J136").Value = 10142.833
$ws.Range("K136").Value = 3340.5651
$ws.Range("L136").Value = 30428.499
$ws.Range("M136").Value = -790.5650999999998
$ws.Range("N136").Value = -35528.499

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4277825
$ws.Range("I122").Value = 4903485.5
$ws.Range("J122").Value = 2481.1667
$ws.Range("K122").Value = 14710456.5
$ws.Range("L122").Value = 7443.500100000001
$ws.Range("M122").Value = -14708006.5
$ws.Range("N122").Value = -12343.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 135747.7
$ws.Range("I132").Value = 210332.53
$ws.Range("J132").Value = 2280.1052
$ws.Range("K132").Value = 630997.59
$ws.Range("L132").Value = 6840.3156
$ws.Range("M132").Value = -628467.59
$ws.Range("N132").Value = -11900.3156

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1268.8928
$ws.Range("I136").Value = 834.9286
$ws.Range("K136").Value = 2504.7858
$ws.Range("M136").Value = 45.21420000000035
